$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.703.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "'3.626.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'583.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "'175.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("D7").Value = "'0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.63%  "

$ws.Range("D8").Value = "'3.622.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "'0.195"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.88%  "

$ws.Range("D11").Value = "'6.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.04%  "

$ws.Range("D12").Value = "'0.618"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").Value = "'48.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("D14").Value = "'0.0000282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").Value = "'677.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "'4.203.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "'9.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "'3.617.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'70.569.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").Value = "'0.123"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "'17.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.34%  "

$ws.Range("D22").Value = "'11.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("D23").Value = "'0.939"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'17.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.08%  "

$ws.Range("D25").Value = "'99.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.96%  "

$ws.Range("D26").Value = "'3.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").Value = "'2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "

$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'9.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "

$ws.Range("D30").Value = "'34.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.48%  "

$ws.Range("D31").Value = "'9.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "'3.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").Value = "'7.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("D34").Value = "'1.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.92%  "

$ws.Range("D35").Value = "'3.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.17%  "

$ws.Range("D36").Value = "'571.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "

$ws.Range("D37").Value = "'11.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("D38").Value = "'0.108"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("D39").Value = "'58.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "

$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.347"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0451"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.14%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'3.537.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "

$ws.Range("D45").Value = "'34.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.75%  "

$ws.Range("D46").Value = "'0.0₃0729"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.96%  "

$ws.Range("D47").Value = "'2.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.91%  "

$ws.Range("D48").Value = "'2.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.61%  "

$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("D50").Value = "'136.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.88%  "

$ws.Range("D51").Value = "'2.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.30%  "
